$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44432
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 24500
$ws.Range("P2").Value = 1633
$ws.Range("D3").Value = 44418
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 26000
$ws.Range("M3").Value = 25500
$ws.Range("P3").Value = 1700
$ws.Range("D4").Value = 44406
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24520
$ws.Range("P4").Value = 1635
$ws.Range("D5").Value = 44343
$ws.Range("J5").Value = 26
$ws.Range("K5").Value = 23000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 23500
$ws.Range("P5").Value = 1567
$ws.Range("D6").Value = 44341
$ws.Range("J6").Value = 36
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 24500
$ws.Range("P6").Value = 1633
$ws.Range("D7").Value = 44421
$ws.Range("J7").Value = 18
$ws.Range("K7").Value = 24000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 24500
$ws.Range("P7").Value = 1633
$ws.Range("D8").Value = 44383
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13480
$ws.Range("P8").Value = 899
$ws.Range("D9").Value = 44460
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 24480
$ws.Range("P9").Value = 1632
$ws.Range("D10").Value = 44425
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 24000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 24520
$ws.Range("P10").Value = 1635
$ws.Range("D11").Value = 44446
$ws.Range("J11").Value = 34
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24500
$ws.Range("P11").Value = 1633
$ws.Range("D12").Value = 44411
$ws.Range("J12").Value = 34
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 26000
$ws.Range("M12").Value = 25500
$ws.Range("P12").Value = 1700
$ws.Range("D13").Value = 44390
$ws.Range("J13").Value = 34
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 24500
$ws.Range("P13").Value = 1633
$ws.Range("D14").Value = 44397
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = 23000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23500
$ws.Range("P14").Value = 1567
$ws.Range("D15").Value = 44351
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 24500
$ws.Range("P15").Value = 1633
$ws.Range("D16").Value = 44349
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24524
$ws.Range("P16").Value = 1635
$ws.Range("D17").Value = 44435
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 24000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 24500
$ws.Range("P17").Value = 1633
$ws.Range("D18").Value = 44336
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24500
$ws.Range("P18").Value = 1633
$ws.Range("D19").Value = 44413
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 24000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 24480
$ws.Range("P19").Value = 1632
$ws.Range("D20").Value = 44329
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 23000
$ws.Range("L20").Value = 23000
$ws.Range("M20").Value = 23000
$ws.Range("P20").Value = 1533
$ws.Range("D21").Value = 44455
$ws.Range("J21").Value = 18
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 24500
$ws.Range("P21").Value = 1633
$ws.Range("D22").Value = 44463
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 24480
$ws.Range("P22").Value = 1632
$ws.Range("D23").Value = 44400
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = 24500
$ws.Range("P23").Value = 1633
$ws.Range("D24").Value = 44442
$ws.Range("J24").Value = 28
$ws.Range("K24").Value = 24000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 24500
$ws.Range("P24").Value = 1633
$ws.Range("D25").Value = 44453
$ws.Range("J25").Value = 25
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 26000
$ws.Range("M25").Value = 25520
$ws.Range("P25").Value = 1701
$ws.Range("D26").Value = 44449
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 24500
$ws.Range("P26").Value = 1633
$ws.Range("D27").Value = 44385
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14480
$ws.Range("P27").Value = 965
$ws.Range("D28").Value = 44428
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25500
$ws.Range("P28").Value = 1700
